# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets, reflecting the latest scrape at
# commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row => new F value
$exhibitionUpdates = @{
    2  = 64
    3  = 364
    6  = 13565
    8  = 64
    9  = 5494
    10 = 566
    12 = 29
    14 = 1217
    16 = 159
    17 = 728
    18 = 2895
    19 = 9146
    20 = 1177
    21 = 3679
    23 = 60
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (all types) - row => new F value
$allTypesUpdates = @{
    2  = 64
    3  = 364
    7  = 13565
    9  = 64
    10 = 5494
    11 = 566
    13 = 29
    15 = 1217
    17 = 159
    18 = 728
    19 = 2895
    21 = 9146
    22 = 1177
    23 = 3679
    25 = 60
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
